# Insert a new data row at row 313 (pushing the existing rows 313:418 down
# to 314:419) and populate it with the new "Apio" price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(313).Insert()

$ws.Cells.Item(313, 1).Value = 11
$ws.Cells.Item(313, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(313, 3).Value = "Bíobío"
$ws.Cells.Item(313, 4).Value = 44985
$ws.Cells.Item(313, 5).Value = 8
$ws.Cells.Item(313, 6).Value = 100112017
$ws.Cells.Item(313, 7).Value = "Apio"
$ws.Cells.Item(313, 8).Value = "Americana (o)"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 300
$ws.Cells.Item(313, 11).Value = 8000
$ws.Cells.Item(313, 12).Value = 8500
$ws.Cells.Item(313, 13).Value = 8250
$ws.Cells.Item(313, 14).Value = "$/docena de matas"
$ws.Cells.Item(313, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(313, 16).Value = 1375
$ws.Cells.Item(313, 17).Value = 6
$ws.Cells.Item(313, 18).Value = "Hortaliza"
